# Commit: "it should create a new group when payer login with a single clone scene"
# Column K on the scene sheet is "Share" — flip it from 1 (TRUE) to 0 (FALSE)
# for every scene row that currently shares a group (except the one that is
# already 0), so a payer who clones a single-instance scene gets put into a
# brand-new group instead of being shared into an existing one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K10").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("K15").Value = 0

# Move the active selection to K10 (matches the saved sheet view state).
$ws.Range("K10").Select() | Out-Null
